$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the two date values from I1/J1 (keep the cell style)
$ws.Range("I1").Value = $null
$ws.Range("J1").Value = $null

# Column width formatting: columns B..O (2-15) become a single uniform width
$ws.Range("B1:O1").EntireColumn.ColumnWidth = 3.625

# Add legend text into column B next to the rows that used to hold it in column A (14-16)
$ws.Range("B11").Value = $ws.Range("A14").Value2
$ws.Range("B12").Value = $ws.Range("A15").Value2
$ws.Range("B13").Value = $ws.Range("A16").Value2

# Remove the old standalone legend rows 14-16
$ws.Range("A14:P16").Delete()

# Update the selection to match the new used range
$ws.Range("A1:P13").Select()
